# Add a new "diesel" fuel row (row 22) to the Fuels sheet, following the
# same layout/formula pattern as the existing charcoal-upstream rows above.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fuels")

$ws.Cells.Item(22, 1).Value = "diesel"               # A22 fuel
$ws.Cells.Item(22, 2).Value = 45.6                    # B22 HHV
$ws.Cells.Item(22, 3).Value = 28.435                  # C22 LHV
$ws.Cells.Item(22, 4).Formula = "=G22*(44/12)"        # D22 CO2
$ws.Cells.Item(22, 5).Value = 0                       # E22 H2O
$ws.Cells.Item(22, 7).Value = 0.86                    # G22 C %
$ws.Cells.Item(22, 12).Value = "ecoinvent 2.2"        # L22 meta-source

# Reflect the author's final on-screen selection after adding the row.
$ws.Range("L24").Select()
